# Update Name of Algo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = 16.4295
$ws.Range("C7").Value = -12.6403
$ws.Range("A9").Value = -21.7488
$ws.Range("C12").Value = -10.7328
$ws.Range("A13").Value = -22.19370000000001
$ws.Range("C14").Value = -13.6705
$ws.Range("E15").Value = 16.59069999999999
$ws.Range("A16").Value = -21.44299999999998
$ws.Range("A18").Value = -22.25190000000001
$ws.Range("C19").Value = -12.09870000000001
$ws.Range("A20").Value = -21.54089999999998
$ws.Range("A26").Value = -21.57199999999998
$ws.Range("C26").Value = -12.569
$ws.Range("A27").Value = -21.57989999999998
$ws.Range("C27").Value = -11.9041
$ws.Range("E28").Value = 16.03949999999999
$ws.Range("A29").Value = -20.66989999999998
$ws.Range("C29").Value = -11.2421
$ws.Range("E33").Value = 17.02770000000002
$ws.Range("A35").Value = -21.37999999999998
$ws.Range("E35").Value = 16.54360000000001
$ws.Range("A36").Value = -20.73769999999998
$ws.Range("C37").Value = -12.44720000000001
$ws.Range("C38").Value = -12.436
$ws.Range("E38").Value = 16.3482
$ws.Range("E43").Value = 17.24410000000002
$ws.Range("E44").Value = 16.7109
$ws.Range("A45").Value = -21.53619999999999
$ws.Range("E45").Value = 16.43609999999999
$ws.Range("C47").Value = -12.45949999999999
$ws.Range("E47").Value = 16.26199999999999
$ws.Range("C51").Value = -11.8419
$ws.Range("E51").Value = 17.18670000000001
$ws.Range("C52").Value = -11.342
$ws.Range("E54").Value = 16.49010000000001
$ws.Range("A55").Value = -22.4634
$ws.Range("C55").Value = -13.64239999999999
$ws.Range("A57").Value = -22.34950000000001
$ws.Range("E57").Value = 16.6968
$ws.Range("E62").Value = 16.45529999999999
$ws.Range("E63").Value = 18.41110000000002
$ws.Range("E67").Value = 17.03750000000002
$ws.Range("A69").Value = -21.7218
$ws.Range("C69").Value = -11.0469
$ws.Range("C70").Value = -12.0214
$ws.Range("E70").Value = 17.42370000000002
$ws.Range("A76").Value = -19.5306
$ws.Range("C76").Value = -13.00500000000001
$ws.Range("A78").Value = -20.23409999999998
$ws.Range("C81").Value = -14.0082
$ws.Range("E81").Value = 16.5392
$ws.Range("A82").Value = -22.2163
$ws.Range("A83").Value = -21.9295
$ws.Range("C83").Value = -13.1138
$ws.Range("E88").Value = 16.31529999999999
$ws.Range("A93").Value = -20.62279999999998
$ws.Range("C94").Value = -10.5349
$ws.Range("E96").Value = 16.20799999999999
$ws.Range("A97").Value = -22.1478
$ws.Range("E99").Value = 16.71690000000001
$ws.Range("C100").Value = -11.6524
$ws.Range("C102").Value = -12.29350000000001
